# Generate Report for Handoff
# Update the "aed6ee9e-a8c2-4f30-9487-fb40fd7a91e0" file row (row 7) on every
# sheet with its newly generated handoff timestamp.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest Handoff Date" column (D) for the aed6ee9e row.
$wsOverview.Range("D7").Value = "2016-30-20 10:30:38"

# zh-cn sheet: "Latest Handoff Datetime" column (E) for the aed6ee9e row.
$wsZhCn.Range("E7").Value = "2016-03-20 10:30:35"

# de-de sheet: "Latest Handoff Datetime" column (E) for the aed6ee9e row.
$wsDeDe.Range("E7").Value = "2016-03-20 10:30:38"
